# Applies the "Elimina antiguos EC y agrega nuevos y modifica Antigua BD" edit:
#  - Updates the "VALOR MORA" total (E11) and the "Cant. Periodos" count (F13)
#  - Inserts a new detail row for worker STEFFANY ANDREA MEZA RENGIFO, period 2509
#    (duplicate of the existing 2508 row, same amounts), pushing the footer rows down

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row right after the last detail row (32) and copy that row's
#    values + formatting (including the bottom-border "last row" style) into it.
$ws.Rows.Item(33).Insert()
$ws.Range("B32:J32").Copy($ws.Range("B33:J33"))

# 2. Row 32 is no longer the last row of the table, so it should go back to the
#    regular interior-row formatting (same as the row above it).
$ws.Range("B31:J31").Copy()
$ws.Range("B32:J32").PasteSpecial(-4122)

# 3. The new row is for period 2509 (same worker/doc/salary/valor mora as 2508).
$ws.Range("E33").Value = "2509"

# 4. Update the summary figures affected by the new period.
$ws.Range("E11").Value = 746080
$ws.Range("F13").Value = 18
